$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '247.34'
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '22.62'
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.286'
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.05731'
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.8101'
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.8664'
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1430'
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07338'
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.03117'
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.09393'
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '3.915'
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.001591'
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.04809'
$cell.Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E17").Value = '16OneONEWorstin24h'
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.0005852'
$cell.Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("E18").Value = '17TigerCashTCH'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.006147'
$cell.Style = "Normal"

# Row 19
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.005143'
$cell.Style = "Normal"

# Row 20
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("E20").Value = '19BitKanKAN'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0009974'
$cell.Style = "Normal"

# Row 21
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("E21").Value = '20NitroExNTX'
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.0001501'
$cell.Style = "Normal"

# Row 22
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("E22").Value = '21LEOLEO'
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '3.730'
$cell.Style = "Normal"

# Row 23
$ws.Range("B23").Value = 'KuCoinToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("E23").Value = '22KuCoinTokenKCS'
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.319'
$cell.Style = "Normal"

# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.185'
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.03938'
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.006766'
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.1067'
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.003201'
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.008148'
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.00005597'
$cell.Style = "Normal"

# Row 47
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.6002'
$cell.Style = "Normal"

# Row 48
$ws.Range("E48").Value = '47BOLOBOLO'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.1769'
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.01010'
$cell.Style = "Normal"
